$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill level")

# The three drop-down (form-control) selectors on the "Skill level" sheet are
# linked to cells B4, B6 and B8. Update the linked cells to the new selected
# values; all dependent formulas (D4:G4, D6:G6, D8:G8, the D15:G15 sums and
# the D17:G17 results) recalculate automatically.
$ws.Range("B4").Value = 2
$ws.Range("B6").Value = 6
$ws.Range("B8").Value = 22

# Keep the combo boxes' reported selection in sync with their linked cells
# (best effort - mirrors what Excel does internally when a linked cell
# driving a Forms combo box changes).
foreach ($pair in @(
        @{ Name = "Vervolgkeuzelijst 4";  Index = 2 },
        @{ Name = "Vervolgkeuzelijst 5";  Index = 6 },
        @{ Name = "Vervolgkeuzelijst 11"; Index = 22 }
    )) {
    try {
        $shp = $ws.Shapes.Item($pair.Name)
        $shp.ControlFormat.ListIndex = $pair.Index
    } catch {
    }
}

$excel.Calculate()
